$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126 (shifts existing rows 126-212 down to 127-213)
$ws.Rows.Item(126).Insert()

# Populate the new row 126 with its data
$ws.Cells.Item(126, 1).Value = 11
$ws.Cells.Item(126, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(126, 3).Value = "Bíobío"
$ws.Cells.Item(126, 4).Value = 45090
$ws.Cells.Item(126, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 5).Value = 8
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100102
$ws.Cells.Item(126, 8).Value = "Cítricos"
$ws.Cells.Item(126, 9).Value = 100102004
$ws.Cells.Item(126, 10).Value = "Mandarina"
$ws.Cells.Item(126, 11).Value = "Clementina"
$ws.Cells.Item(126, 12).Value = "Primera"
$ws.Cells.Item(126, 13).Value = 190
$ws.Cells.Item(126, 14).Value = 8000
$ws.Cells.Item(126, 15).Value = 9000
$ws.Cells.Item(126, 16).Value = 8526
$ws.Cells.Item(126, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(126, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(126, 19).Value = 853
$ws.Cells.Item(126, 20).Value = 10
